$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743717012329/Ggg0CQFNgb-camera-video-segment-driver-1743717014829.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMKBZOCVG%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T160310Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCJq%2Fut6d%2FMNBoCh9xrwHFf5AIs9EVH1IAJygqXASBnTQIgDY17Xns4w1eL7XLrvvB%2BpHKlgc3FX2Wdpj8WlMFJEZkq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDMxYMZ435Ef20k%2BE7Sq6A0EqpuvD353ulDczjQ2%2BUDfQyY%2Bqm8Lio%2FmW%2FQl%2FvafFtl1xElz3beZAk9pVOcb0sqhuzlXGdxa58IPlfjtUNi9WbvNb1gjqPbSsP%2BChY6Th3LhDJa6z%2FaqA%2FE4rY8ycy7qCLLxoV45FaFjjUablyMTV7M36HT2ZK1g%2F8dTUhvAsR%2FGI8azi5RQh0HwODjOgfQRBzIC5uLidn1RqSzw6agT02yI5jWuoCjtvMmZTN%2FyCyjYEMDBSM3%2BBBwVTpnDqLZ9J94Xm4v4yjnqXPyn%2Bx%2BInykX5r0KXv%2B3eiFD%2BVjn0jZ3l1cvp1k2%2BAC%2BWOzfWJrUNFQPKIU%2B%2F2gVeeeubq4oi9Fh2418gYlOoGXdCzTMn06X2rg6THZnNsHaB69E%2Fbv4crN%2F%2BxCOf7em5UI2WRGJNJAulBgM9Sjn%2FjiM0dK2DPYr7uFPoOOA4yrSXLPKxnV0FfC3cIzDwO8i8A%2FDQTWe%2BsiQZaZnXyUzooVHIxfwJqLvQNlOLNzhMh%2F9JucT2VMtOT%2BfIr3Qa%2ByDxPNvJ%2FD23NPmUUWxqcZnprbyratA7IgwMX13EJKviEE3xJQQCqtGexWBShVTmHm4wzKW%2FvwY6pQHGW2qMxH04W66L797lnvlnBoZZP7cmN4ATzl2OTI8PAtyVvXwdIuJSBDgfH0FtC3zki762wNVlJodjZb%2Fx%2FoE0ntaws6YBUAD%2BfrQBSAinTYErGMkQsHCEoVAjm1u9exjZJZROqlEigrj3kyGYTfTMJYnNC%2FM%2FSzLQ3LNUQKNl4Lhvv2AQkF5eugmzywZ%2F4iaDoli6%2BKN%2BD0ma71tq60X8L0bkoKA%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2000%3A03%3A10%20GMT&X-Amz-Signature=636c570c59af025ac471b9e5f9d62a12f93af626ade10f9a0f9350f298aaf04b'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743713574485/FgwSeOefdQ-camera-video-segment-driver-1743713576985.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMKBZOCVG%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T160310Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCJq%2Fut6d%2FMNBoCh9xrwHFf5AIs9EVH1IAJygqXASBnTQIgDY17Xns4w1eL7XLrvvB%2BpHKlgc3FX2Wdpj8WlMFJEZkq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDMxYMZ435Ef20k%2BE7Sq6A0EqpuvD353ulDczjQ2%2BUDfQyY%2Bqm8Lio%2FmW%2FQl%2FvafFtl1xElz3beZAk9pVOcb0sqhuzlXGdxa58IPlfjtUNi9WbvNb1gjqPbSsP%2BChY6Th3LhDJa6z%2FaqA%2FE4rY8ycy7qCLLxoV45FaFjjUablyMTV7M36HT2ZK1g%2F8dTUhvAsR%2FGI8azi5RQh0HwODjOgfQRBzIC5uLidn1RqSzw6agT02yI5jWuoCjtvMmZTN%2FyCyjYEMDBSM3%2BBBwVTpnDqLZ9J94Xm4v4yjnqXPyn%2Bx%2BInykX5r0KXv%2B3eiFD%2BVjn0jZ3l1cvp1k2%2BAC%2BWOzfWJrUNFQPKIU%2B%2F2gVeeeubq4oi9Fh2418gYlOoGXdCzTMn06X2rg6THZnNsHaB69E%2Fbv4crN%2F%2BxCOf7em5UI2WRGJNJAulBgM9Sjn%2FjiM0dK2DPYr7uFPoOOA4yrSXLPKxnV0FfC3cIzDwO8i8A%2FDQTWe%2BsiQZaZnXyUzooVHIxfwJqLvQNlOLNzhMh%2F9JucT2VMtOT%2BfIr3Qa%2ByDxPNvJ%2FD23NPmUUWxqcZnprbyratA7IgwMX13EJKviEE3xJQQCqtGexWBShVTmHm4wzKW%2FvwY6pQHGW2qMxH04W66L797lnvlnBoZZP7cmN4ATzl2OTI8PAtyVvXwdIuJSBDgfH0FtC3zki762wNVlJodjZb%2Fx%2FoE0ntaws6YBUAD%2BfrQBSAinTYErGMkQsHCEoVAjm1u9exjZJZROqlEigrj3kyGYTfTMJYnNC%2FM%2FSzLQ3LNUQKNl4Lhvv2AQkF5eugmzywZ%2F4iaDoli6%2BKN%2BD0ma71tq60X8L0bkoKA%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2000%3A03%3A10%20GMT&X-Amz-Signature=0976ddd7b4dc99c48d6ebc330fdbab02f4366e38a06a961c8f886281dc045cb0'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712999664/UzzO0YQTJX-camera-video-segment-driver-1743713002164.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMKBZOCVG%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T160310Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCJq%2Fut6d%2FMNBoCh9xrwHFf5AIs9EVH1IAJygqXASBnTQIgDY17Xns4w1eL7XLrvvB%2BpHKlgc3FX2Wdpj8WlMFJEZkq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDMxYMZ435Ef20k%2BE7Sq6A0EqpuvD353ulDczjQ2%2BUDfQyY%2Bqm8Lio%2FmW%2FQl%2FvafFtl1xElz3beZAk9pVOcb0sqhuzlXGdxa58IPlfjtUNi9WbvNb1gjqPbSsP%2BChY6Th3LhDJa6z%2FaqA%2FE4rY8ycy7qCLLxoV45FaFjjUablyMTV7M36HT2ZK1g%2F8dTUhvAsR%2FGI8azi5RQh0HwODjOgfQRBzIC5uLidn1RqSzw6agT02yI5jWuoCjtvMmZTN%2FyCyjYEMDBSM3%2BBBwVTpnDqLZ9J94Xm4v4yjnqXPyn%2Bx%2BInykX5r0KXv%2B3eiFD%2BVjn0jZ3l1cvp1k2%2BAC%2BWOzfWJrUNFQPKIU%2B%2F2gVeeeubq4oi9Fh2418gYlOoGXdCzTMn06X2rg6THZnNsHaB69E%2Fbv4crN%2F%2BxCOf7em5UI2WRGJNJAulBgM9Sjn%2FjiM0dK2DPYr7uFPoOOA4yrSXLPKxnV0FfC3cIzDwO8i8A%2FDQTWe%2BsiQZaZnXyUzooVHIxfwJqLvQNlOLNzhMh%2F9JucT2VMtOT%2BfIr3Qa%2ByDxPNvJ%2FD23NPmUUWxqcZnprbyratA7IgwMX13EJKviEE3xJQQCqtGexWBShVTmHm4wzKW%2FvwY6pQHGW2qMxH04W66L797lnvlnBoZZP7cmN4ATzl2OTI8PAtyVvXwdIuJSBDgfH0FtC3zki762wNVlJodjZb%2Fx%2FoE0ntaws6YBUAD%2BfrQBSAinTYErGMkQsHCEoVAjm1u9exjZJZROqlEigrj3kyGYTfTMJYnNC%2FM%2FSzLQ3LNUQKNl4Lhvv2AQkF5eugmzywZ%2F4iaDoli6%2BKN%2BD0ma71tq60X8L0bkoKA%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2000%3A03%3A10%20GMT&X-Amz-Signature=ef372ed673f5b52ef41ae0778e2f36e0e636a8261f54289e5a958ac26531b9a6'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712546260/Jh2p7rn04n-camera-video-segment-driver-1743712548760.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMKBZOCVG%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T160310Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCJq%2Fut6d%2FMNBoCh9xrwHFf5AIs9EVH1IAJygqXASBnTQIgDY17Xns4w1eL7XLrvvB%2BpHKlgc3FX2Wdpj8WlMFJEZkq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDMxYMZ435Ef20k%2BE7Sq6A0EqpuvD353ulDczjQ2%2BUDfQyY%2Bqm8Lio%2FmW%2FQl%2FvafFtl1xElz3beZAk9pVOcb0sqhuzlXGdxa58IPlfjtUNi9WbvNb1gjqPbSsP%2BChY6Th3LhDJa6z%2FaqA%2FE4rY8ycy7qCLLxoV45FaFjjUablyMTV7M36HT2ZK1g%2F8dTUhvAsR%2FGI8azi5RQh0HwODjOgfQRBzIC5uLidn1RqSzw6agT02yI5jWuoCjtvMmZTN%2FyCyjYEMDBSM3%2BBBwVTpnDqLZ9J94Xm4v4yjnqXPyn%2Bx%2BInykX5r0KXv%2B3eiFD%2BVjn0jZ3l1cvp1k2%2BAC%2BWOzfWJrUNFQPKIU%2B%2F2gVeeeubq4oi9Fh2418gYlOoGXdCzTMn06X2rg6THZnNsHaB69E%2Fbv4crN%2F%2BxCOf7em5UI2WRGJNJAulBgM9Sjn%2FjiM0dK2DPYr7uFPoOOA4yrSXLPKxnV0FfC3cIzDwO8i8A%2FDQTWe%2BsiQZaZnXyUzooVHIxfwJqLvQNlOLNzhMh%2F9JucT2VMtOT%2BfIr3Qa%2ByDxPNvJ%2FD23NPmUUWxqcZnprbyratA7IgwMX13EJKviEE3xJQQCqtGexWBShVTmHm4wzKW%2FvwY6pQHGW2qMxH04W66L797lnvlnBoZZP7cmN4ATzl2OTI8PAtyVvXwdIuJSBDgfH0FtC3zki762wNVlJodjZb%2Fx%2FoE0ntaws6YBUAD%2BfrQBSAinTYErGMkQsHCEoVAjm1u9exjZJZROqlEigrj3kyGYTfTMJYnNC%2FM%2FSzLQ3LNUQKNl4Lhvv2AQkF5eugmzywZ%2F4iaDoli6%2BKN%2BD0ma71tq60X8L0bkoKA%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2000%3A03%3A10%20GMT&X-Amz-Signature=c21ff43c7ac87245b98790689425dd8eb10997799508829bc4bc4c0f0472cc27'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743703000742/spTbAbfADq-camera-video-segment-driver-1743703003242.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMKBZOCVG%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T160310Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCJq%2Fut6d%2FMNBoCh9xrwHFf5AIs9EVH1IAJygqXASBnTQIgDY17Xns4w1eL7XLrvvB%2BpHKlgc3FX2Wdpj8WlMFJEZkq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDMxYMZ435Ef20k%2BE7Sq6A0EqpuvD353ulDczjQ2%2BUDfQyY%2Bqm8Lio%2FmW%2FQl%2FvafFtl1xElz3beZAk9pVOcb0sqhuzlXGdxa58IPlfjtUNi9WbvNb1gjqPbSsP%2BChY6Th3LhDJa6z%2FaqA%2FE4rY8ycy7qCLLxoV45FaFjjUablyMTV7M36HT2ZK1g%2F8dTUhvAsR%2FGI8azi5RQh0HwODjOgfQRBzIC5uLidn1RqSzw6agT02yI5jWuoCjtvMmZTN%2FyCyjYEMDBSM3%2BBBwVTpnDqLZ9J94Xm4v4yjnqXPyn%2Bx%2BInykX5r0KXv%2B3eiFD%2BVjn0jZ3l1cvp1k2%2BAC%2BWOzfWJrUNFQPKIU%2B%2F2gVeeeubq4oi9Fh2418gYlOoGXdCzTMn06X2rg6THZnNsHaB69E%2Fbv4crN%2F%2BxCOf7em5UI2WRGJNJAulBgM9Sjn%2FjiM0dK2DPYr7uFPoOOA4yrSXLPKxnV0FfC3cIzDwO8i8A%2FDQTWe%2BsiQZaZnXyUzooVHIxfwJqLvQNlOLNzhMh%2F9JucT2VMtOT%2BfIr3Qa%2ByDxPNvJ%2FD23NPmUUWxqcZnprbyratA7IgwMX13EJKviEE3xJQQCqtGexWBShVTmHm4wzKW%2FvwY6pQHGW2qMxH04W66L797lnvlnBoZZP7cmN4ATzl2OTI8PAtyVvXwdIuJSBDgfH0FtC3zki762wNVlJodjZb%2Fx%2FoE0ntaws6YBUAD%2BfrQBSAinTYErGMkQsHCEoVAjm1u9exjZJZROqlEigrj3kyGYTfTMJYnNC%2FM%2FSzLQ3LNUQKNl4Lhvv2AQkF5eugmzywZ%2F4iaDoli6%2BKN%2BD0ma71tq60X8L0bkoKA%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2000%3A03%3A10%20GMT&X-Amz-Signature=90cb2d80efc9cf728f4c100f8861215e2aa1356da1adcbbb669e665c32efeb44'

$wb.Save()